$d = $word.ActiveDocument

function Mark-Checkbox($searchText) {
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $searchText"
        return
    }
    $start = $r.Start
    # The underscore immediately after the leading one is turned into an "X"
    # (i.e. "___" -> "_X_"), splitting the run into three runs just like
    # Word does when typing over a mid-run selection.
    $charRange = $d.Range($start + 1, $start + 2)
    $charRange.Font.Bold = $true
    $charRange.Text = "X"
    $charRange2 = $d.Range($start + 1, $start + 2)
    $charRange2.Font.Bold = $false
}

Mark-Checkbox("___ Uploaded    ___ N/A")
Mark-Checkbox("___ Completed     ___ N/A")

Write-Host "done"
